# Logged Week 16 and performed season sim from Week 17
#
# Rushing sheet: bump cumulative rushing totals for Week 16 logging.
# Receiving sheet: move J.Waddle's row into its normal sorted slot (was parked
# at the bottom of the table) with updated totals, add newly-logged player
# T.Lewis, and bump several other players' cumulative receiving totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# M.Gaskin (row 4) - 2DATT
$rushing.Range("D4").Value = 62

# D.Johnson (row 8)
$rushing.Range("C8").Value = 21
$rushing.Range("D8").Value = 14
$rushing.Range("E8").Value = 4
$rushing.Range("F8").Value = 7

# P.Lindsay (row 9)
$rushing.Range("C9").Value = 13
$rushing.Range("D9").Value = 6
$rushing.Range("E9").Value = 4
$rushing.Range("F9").Value = 3

# P.Williams (row 11)
$rushing.Range("D11").Value = 2
$rushing.Range("E11").Value = 1

# ---------------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# J.Waddle's row had been parked at the very bottom of the table (row 18);
# pull it out and re-insert it in its normal sorted position (row 8, right
# after D.Parker) with the Week-16-updated totals.
$receiving.Rows.Item(18).Delete()

$receiving.Rows.Item(8).Insert()
$receiving.Range("A8").Value = 6
$receiving.Range("B8").Value = "J.Waddle"
$receiving.Range("C8").Value = 118
$receiving.Range("D8").Value = 97
$receiving.Range("E8").Value = 21
$receiving.Range("F8").Value = 10
$receiving.Range("G8").Value = 16
$receiving.Range("H8").Value = 13
$receiving.Range("A8:H8").Font.Bold = $true
$receiving.Range("A8:H8").HorizontalAlignment = -4108
$receiving.Range("A8:H8").VerticalAlignment = -4160

# M.Gaskin (row 2)
$receiving.Range("C2").Value = 57
$receiving.Range("D2").Value = 44
$receiving.Range("G2").Value = 9

# M.Hollins (now row 12, after the J.Waddle insert)
$receiving.Range("E12").Value = 9
$receiving.Range("F12").Value = 4

# Newly logged player T.Lewis, inserted right after K.Merritt (row 14) and
# before M.Gesicki (row 15).
$receiving.Rows.Item(15).Insert()
$receiving.Range("A15").Value = 13
$receiving.Range("B15").Value = "T.Lewis"
$receiving.Range("C15").Value = 2
$receiving.Range("D15").Value = 1
$receiving.Range("E15").Value = 0
$receiving.Range("F15").Value = 0
$receiving.Range("G15").Value = 1
$receiving.Range("H15").Value = 0
$receiving.Range("A15:H15").Font.Bold = $true
$receiving.Range("A15:H15").HorizontalAlignment = -4108
$receiving.Range("A15:H15").VerticalAlignment = -4160

# M.Gesicki (now row 16)
$receiving.Range("C16").Value = 81
$receiving.Range("D16").Value = 55
$receiving.Range("E16").Value = 21
$receiving.Range("G16").Value = 10
$receiving.Range("H16").Value = 8

# D.Smythe (now row 19, last row)
$receiving.Range("C19").Value = 30
$receiving.Range("D19").Value = 25
